$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 17072.182
$ws.Range("I6").Value = 18199.5
$ws.Range("K6").Value = 54598.5
$ws.Range("M6").Value = -54486.5
$ws.Range("H19").Value = 909.5454999999999
$ws.Range("I19").Value = 897.8
$ws.Range("K19").Value = 897.8
$ws.Range("M19").Value = -722.8
$ws.Range("H113").Value = 2002.5
$ws.Range("I113").Value = 2002.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2002.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1251.5
$ws.Range("N113").ClearContents()
$ws.Range("H118").Value = 617.3333
$ws.Range("I118").Value = 677.5714
$ws.Range("J118").Value = 406.5
$ws.Range("K118").Value = 2032.7142
$ws.Range("L118").Value = 1219.5
$ws.Range("M118").Value = -375.7142000000001
$ws.Range("N118").Value = -4533.5
$ws.Range("H132").Value = 2164.3333
$ws.Range("I132").Value = 2340.7778
$ws.Range("J132").Value = 1635
$ws.Range("K132").Value = 7022.3334
$ws.Range("L132").Value = 4905
$ws.Range("M132").Value = -4492.3334
$ws.Range("N132").Value = -9965
$ws.Range("H137").Value = 12673.037
$ws.Range("I137").Value = 2075.4375
$ws.Range("K137").Value = 6226.3125
$ws.Range("M137").Value = -3676.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2711.1875
$ws.Range("I74").Value = 1071.3182
$ws.Range("K74").Value = 1071.3182
$ws.Range("M74").Value = -197.3181999999999
$ws.Range("H77").Value = 2711.1875
$ws.Range("I77").Value = 1071.3182
$ws.Range("K77").Value = 5356.590999999999
$ws.Range("M77").Value = -988.5909999999994

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 491.5
$ws.Range("I22").Value = 433
$ws.Range("K22").Value = 433
$ws.Range("M22").Value = -260
$ws.Range("H134").Value = 7697.3667
$ws.Range("I134").Value = 5020.3335
$ws.Range("J134").Value = 13943.777
$ws.Range("K134").Value = 15061.0005
$ws.Range("L134").Value = 41831.331
$ws.Range("M134").Value = -12526.0005
$ws.Range("N134").Value = -46901.331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3385.9412
$ws.Range("J31").Value = 3843.3076
$ws.Range("L31").Value = 3843.3076
$ws.Range("N31").Value = -4433.3076
$ws.Range("H34").Value = 3385.9412
$ws.Range("J34").Value = 3843.3076
$ws.Range("L34").Value = 3843.3076
$ws.Range("N34").Value = -4247.3076
$ws.Range("H70").Value = 99889
$ws.Range("J70").Value = 99889
$ws.Range("L70").Value = 99889
$ws.Range("N70").Value = -100519
$ws.Range("H73").Value = 99889
$ws.Range("J73").Value = 99889
$ws.Range("L73").Value = 99889
$ws.Range("N73").Value = -102073
$ws.Range("H80").Value = 69989
$ws.Range("J80").Value = 69989
$ws.Range("L80").Value = 69989
$ws.Range("N80").Value = -72235
$ws.Range("H83").Value = 69989
$ws.Range("J83").Value = 69989
$ws.Range("L83").Value = 209967
$ws.Range("N83").Value = -221199
$ws.Range("H97").Value = 54397
$ws.Range("J97").Value = 19500
$ws.Range("L97").Value = 19500
$ws.Range("N97").Value = -21482
$ws.Range("H102").Value = 99848
$ws.Range("J102").Value = 99848
$ws.Range("L102").Value = 99848
$ws.Range("N102").Value = -104716
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -55242
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H115").Value = 75000
$ws.Range("J115").Value = 75000
$ws.Range("L115").Value = 75000
$ws.Range("N115").Value = -77350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1591
$ws.Range("J131").Value = 1790.9615
$ws.Range("L131").Value = 5372.8845
$ws.Range("N131").Value = -15452.8845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12460.333
$ws.Range("I43").Value = 4953.4
$ws.Range("J43").Value = 49995
$ws.Range("K43").Value = 4953.4
$ws.Range("L43").Value = 49995
$ws.Range("M43").Value = -4802.4
$ws.Range("N43").Value = -50297
$ws.Range("H57").Value = 32569.363
$ws.Range("I57").Value = 4055
$ws.Range("J57").Value = 38905.89
$ws.Range("K57").Value = 4055
$ws.Range("L57").Value = 38905.89
$ws.Range("M57").Value = -3235
$ws.Range("N57").Value = -40545.89
$ws.Range("H58").Value = 19477.375
$ws.Range("I58").Value = 12409.75
$ws.Range("J58").Value = 26545
$ws.Range("K58").Value = 12409.75
$ws.Range("L58").Value = 26545
$ws.Range("M58").Value = -12132.75
$ws.Range("N58").Value = -27099
$ws.Range("H70").Value = 7548.846
$ws.Range("I70").Value = 6089.7144
$ws.Range("K70").Value = 6089.7144
$ws.Range("M70").Value = -5819.7144
$ws.Range("H73").Value = 7548.846
$ws.Range("I73").Value = 6089.7144
$ws.Range("K73").Value = 6089.7144
$ws.Range("M73").Value = -5153.7144
$ws.Range("H80").Value = 4308.1113
$ws.Range("I80").Value = 3995.4
$ws.Range("K80").Value = 3995.4
$ws.Range("M80").Value = -2997.4
$ws.Range("H83").Value = 4308.1113
$ws.Range("I83").Value = 3995.4
$ws.Range("K83").Value = 19977
$ws.Range("M83").Value = -14985
$ws.Range("H122").Value = 1088.6666
$ws.Range("I122").Value = 833.25
$ws.Range("J122").Value = 1599.5
$ws.Range("K122").Value = 2499.75
$ws.Range("L122").Value = 4798.5
$ws.Range("M122").Value = -49.75
$ws.Range("N122").Value = -9698.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1922.3334
$ws.Range("I16").Value = 2303.5
$ws.Range("K16").Value = 2303.5
$ws.Range("M16").Value = -2133.5
$ws.Range("H22").Value = 2820.718
$ws.Range("I22").Value = 2409.6365
$ws.Range("J22").Value = 5081.6665
$ws.Range("K22").Value = 2409.6365
$ws.Range("L22").Value = 5081.6665
$ws.Range("M22").Value = -2114.6365
$ws.Range("N22").Value = -5671.6665
$ws.Range("H27").Value = 2820.718
$ws.Range("I27").Value = 2409.6365
$ws.Range("J27").Value = 5081.6665
$ws.Range("K27").Value = 2409.6365
$ws.Range("L27").Value = 5081.6665
$ws.Range("M27").Value = -2302.6365
$ws.Range("N27").Value = -5295.6665
$ws.Range("H69").Value = 99999
$ws.Range("J69").Value = 99999
$ws.Range("L69").Value = 99999
$ws.Range("N69").Value = -101621
$ws.Range("H72").Value = 99999
$ws.Range("J72").Value = 99999
$ws.Range("L72").Value = 299997
$ws.Range("N72").Value = -308109
$ws.Range("H76").Value = 19999
$ws.Range("J76").Value = 19999
$ws.Range("L76").Value = 19999
$ws.Range("N76").Value = -20675
$ws.Range("H79").Value = 19999
$ws.Range("J79").Value = 19999
$ws.Range("L79").Value = 19999
$ws.Range("N79").Value = -22339
$ws.Range("H93").Value = 1660.7576
$ws.Range("I93").Value = 1662.7241
$ws.Range("J93").Value = 1646.5
$ws.Range("K93").Value = 1662.7241
$ws.Range("L93").Value = 1646.5
$ws.Range("M93").Value = -414.7240999999999
$ws.Range("N93").Value = -4142.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H136").Value = 11113573
$ws.Range("I136").Value = 11113573
$ws.Range("K136").Value = 33340719
$ws.Range("M136").Value = -33338169

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 687
$ws.Range("I113").Value = 666
$ws.Range("K113").Value = 1998
$ws.Range("M113").Value = 172
$ws.Range("H122").Value = 4880.5454
$ws.Range("J122").Value = 5555
$ws.Range("L122").Value = 16665
$ws.Range("N122").Value = -21565
$ws.Range("H132").Value = 14091.93
$ws.Range("I132").Value = 9718.565000000001
$ws.Range("K132").Value = 29155.695
$ws.Range("M132").Value = -26625.695
$ws.Range("H136").Value = 9093063
$ws.Range("I136").Value = 9093063
$ws.Range("K136").Value = 27279189
$ws.Range("M136").Value = -27276639
